$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.920.62'
$ws.Range("E2").Value = '  -0.17%  '
$ws.Range("D3").Value = '1.669.97'
$ws.Range("E3").Value = '  +0.95%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '214.72'
$ws.Range("E5").Value = '  -0.18%  '
$ws.Range("D6").Value = '0.517'
$ws.Range("E6").Value = '  +0.88%  '
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("E8").Value = '  -0.30%  '
$ws.Range("E9").Value = '  +0.83%  '
$ws.Range("D10").Value = '20.28'
$ws.Range("E10").Value = '  +0.31%  '
$ws.Range("D11").Value = '0.0889'
$ws.Range("E11").Value = '  +1.71%  '
$ws.Range("D12").Value = '1.904.31'
$ws.Range("E12").Value = '  +0.86%  '
$ws.Range("D13").Value = '1.644.32'
$ws.Range("E13").Value = '  -0.72%  '
$ws.Range("E14").Value = '  +0.02%  '
$ws.Range("E15").Value = '  +1.35%  '
$ws.Range("D16").Value = '65.47'
$ws.Range("E16").Value = '  +0.56%  '
$ws.Range("D17").Value = '26.909.02'
$ws.Range("E17").Value = '  -0.24%  '
$ws.Range("D18").Value = '233.29'
$ws.Range("E18").Value = '  -1.07%  '
$ws.Range("D19").Value = '8.01'
$ws.Range("E19").Value = '  +3.89%  '
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("E21").Value = '  -0.06%  '
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("E23").Value = '  -1.77%  '
$ws.Range("D25").Value = '146.24'
$ws.Range("D26").Value = '7.11'
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("D27").Value = '15.93'
$ws.Range("E27").Value = '  +0.80%  '
$ws.Range("E28").Value = '  -0.14%  '
$ws.Range("E29").Value = '  -1.85%  '
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("E31").Value = '  -0.02%  '
$ws.Range("E32").Value = '  +0.17%  '
$ws.Range("D33").Value = '1.457.07'
$ws.Range("E33").Value = '  -5.88%  '
$ws.Range("E34").Value = '  +1.43%  '
$ws.Range("D35").Value = '1.64'
$ws.Range("E35").Value = '  +2.49%  '
$ws.Range("E36").Value = '  +0.01%  '
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("D38").Value = '0.899'
$ws.Range("E38").Value = '  +0.83%  '
$ws.Range("D40").Value = '1.04'
$ws.Range("E40").Value = '  +13.83%  '
$ws.Range("E41").Value = '  -4.43%  '
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("E43").Value = '  +2.22%  '
$ws.Range("D44").Value = '66.22'
$ws.Range("E44").Value = '  +1.23%  '
$ws.Range("D45").Value = '1.811.03'
$ws.Range("E45").Value = '  +0.87%  '
$ws.Range("E46").Value = '  +0.54%  '
$ws.Range("D47").Value = '90.73'
$ws.Range("E47").Value = '  +0.56%  '
$ws.Range("E48").Value = '  +1.41%  '
$ws.Range("D49").Value = '0.101'
$ws.Range("E49").Value = '  +2.78%  '
$ws.Range("E50").Value = '  +0.25%  '
$ws.Range("E51").Value = '  -0.29%  '
